$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new party_normalized (column J) values for the rows that were
# researched/normalized in this pass. Shared strings are created in the same
# order Excel would create them so new parties line up with the source order.

$ws.Range("J2").Value = "Praja Socialist Party"
$ws.Range("J7").Value = "Praja Socialist Party"
$ws.Range("J8").Value = "Praja Socialist Party"
$ws.Range("J14").Value = "Praja Socialist Party"
$ws.Range("J15").Value = "Praja Socialist Party"
$ws.Range("J20").Value = "Praja Socialist Party"
$ws.Range("J25").Value = "Praja Socialist Party"
$ws.Range("J27").Value = "Praja Socialist Party"
$ws.Range("J30").Value = "Praja Socialist Party"
$ws.Range("J3").Value = "Commmunist Party of India"
$ws.Range("J4").Value = "Commmunist Party of India"
$ws.Range("J18").Value = "Commmunist Party of India"
$ws.Range("J32").Value = "Commmunist Party of India"
$ws.Range("J5").Value = "Ganatantra Parishad"
$ws.Range("J6").Value = "Ganatantra Parishad"
$ws.Range("J23").Value = "Ganatantra Parishad"
$ws.Range("J9").Value = "Peasants and Workers Party"
$ws.Range("J10").Value = "Peasants and Workers Party"
$ws.Range("J35").Value = "Peasants and Workers Party"
$ws.Range("J49").Value = "Peasants and Workers Party"
$ws.Range("J12").Value = "???"
$ws.Range("J86").Value = "???"
$ws.Range("J108").Value = "Shiromani Akali Dal"
$ws.Range("J109").Value = "Shiromani Akali Dal"
$ws.Range("J112").Value = "Shiromani Akali Dal"
$ws.Range("J13").Value = "Socialist Party"
$ws.Range("J24").Value = "Socialist Party"
$ws.Range("J28").Value = "Socialist Party"
$ws.Range("J33").Value = "Socialist Party"
$ws.Range("J38").Value = "Socialist Party"
$ws.Range("J46").Value = "Socialist Party"
$ws.Range("J51").Value = "Socialist Party"
$ws.Range("J16").Value = "??"
$ws.Range("J26").Value = "??"
$ws.Range("J19").Value = "Republican Party of India"
$ws.Range("J21").Value = "Republican Party of India"
$ws.Range("J22").Value = "Republican Party of India"
$ws.Range("J36").Value = "Republican Party of India"
$ws.Range("J37").Value = "Republican Party of India"
$ws.Range("J117").Value = "Republican Party of India"
$ws.Range("J29").Value = "Kisan Mazdoor Praja Party"
$ws.Range("J34").Value = "Progressive Independent Party"
$ws.Range("J39").Value = "Jan Sangh"
$ws.Range("J11").Value = "PDF"
$ws.Range("J17").Value = "PDF"
$ws.Range("J31").Value = "Janata Party"
$ws.Range("J42").Value = "Janata Party"
$ws.Range("J43").Value = "Janata Party"
$ws.Range("J44").Value = "Janata Party"
$ws.Range("J45").Value = "Bhartiya Janata Party"
$ws.Range("J50").Value = "Janata Party"
$ws.Range("J52").Value = "Bhartiya Janata Party"
$ws.Range("J53").Value = "Janata Party"
$ws.Range("J56").Value = "Janata Party"
$ws.Range("J57").Value = "Janata Party"
$ws.Range("J58").Value = "Bhartiya Janata Party"
$ws.Range("J59").Value = "Janata Party"
$ws.Range("J61").Value = "Bhartiya Janata Party"
$ws.Range("J62").Value = "Bhartiya Janata Party"
$ws.Range("J63").Value = "Janata Party"
$ws.Range("J65").Value = "Bhartiya Janata Party"
$ws.Range("J67").Value = "Janata Party"
$ws.Range("J68").Value = "Janata Party"
$ws.Range("J69").Value = "Bhartiya Janata Party"
$ws.Range("J70").Value = "Janata Party"
$ws.Range("J71").Value = "Janata Party"
$ws.Range("J72").Value = "Janata Party"
$ws.Range("J73").Value = "Janata Party"
$ws.Range("J74").Value = "Janata Party"
$ws.Range("J75").Value = "Janata Party"
$ws.Range("J76").Value = "Janata Party"
$ws.Range("J77").Value = "Janata Party"
$ws.Range("J78").Value = "Janata Party"
$ws.Range("J79").Value = "Janata Party"
$ws.Range("J80").Value = "Janata Party"
$ws.Range("J81").Value = "Janata Party"
$ws.Range("J82").Value = "Janata Party"
$ws.Range("J83").Value = "Janata Party"
$ws.Range("J84").Value = "Janata Party"
$ws.Range("J85").Value = "Janata Party"
$ws.Range("J87").Value = "Janata Party"
$ws.Range("J88").Value = "Janata Party"
$ws.Range("J89").Value = "Janata Party"
$ws.Range("J90").Value = "Janata Party"
$ws.Range("J91").Value = "Janata Party"
$ws.Range("J92").Value = "Janata Party"
$ws.Range("J93").Value = "Janata Party"
$ws.Range("J94").Value = "Janata Party"
$ws.Range("J95").Value = "Janata Party"
$ws.Range("J96").Value = "Janata Party"
$ws.Range("J97").Value = "Janata Party"
$ws.Range("J98").Value = "Janata Party"
$ws.Range("J100").Value = "Janata Party"
$ws.Range("J101").Value = "Janata Party"
$ws.Range("J103").Value = "Janata Party"

# Row 12 ("SAD" in 1952) could not be confidently normalized -> mark it and
# highlight it in light red so it stands out for a follow-up pass.
$ws.Range("J12").Interior.Color = 14408946

# Row 29 (KMPP) was typed in manually with a different font/size, matching the
# source workbook, which also bumped that row a touch taller.
$ws.Range("J29").Font.Name = "Verdana"
$ws.Range("J29").Font.Size = 12
$ws.Range("J29").Font.Color = 0
$ws.Rows.Item(29).RowHeight = 16

# Turn off the autofilter (left over from the original export) and set the
# print setup like the final workbook.
$ws.AutoFilterMode = $false
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
